$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the target sentence in the last paragraph.
# ------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Because of COVID-19, this project can fuck off ")
if (-not $found) {
    throw "Could not find target sentence"
}
$sentenceStart = $find.Start
$sentenceEnd = $find.End

$becauseOfStart = $sentenceStart
$becauseOfEnd = $sentenceStart + 10     # length of "Because of"

# ------------------------------------------------------------------
# 2. Insert " Thanks to" right after "Because of". At this point it
#    will merge into one run with "Because of" (same formatting) -
#    that's fixed up in step 4 below.
# ------------------------------------------------------------------
$insertionPoint = $d.Range($becauseOfEnd, $becauseOfEnd)
$insertionPoint.InsertAfter(" Thanks to")

$thanksToEnd = $becauseOfEnd + 10       # length of " Thanks to"

# ------------------------------------------------------------------
# 3. Force a run boundary between " Thanks to" and
#    " COVID-19, this project can fuck off " by briefly dropping a
#    bookmark at the seam and removing it again. The bookmark
#    start/end elements physically separate the two <w:r> runs, so
#    once it is deleted the split survives without the usual
#    identical-formatting run merge.
# ------------------------------------------------------------------
$seam = $d.Range($thanksToEnd, $thanksToEnd)
$tempBookmark = $d.Bookmarks.Add("ZZTempSeam", $seam)
$d.Bookmarks("ZZTempSeam").Delete()

# ------------------------------------------------------------------
# 4. Apply strikethrough to just "Because of" - this splits it off
#    from " Thanks to" into its own run (differing formatting means
#    no merge occurs).
# ------------------------------------------------------------------
$becauseOfRange = $d.Range($becauseOfStart, $becauseOfEnd)
$becauseOfRange.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 5. Move the _GoBack bookmark from the end of the document (after
#    the emoji run) to right after "it's" / before the gramEnd
#    proofing mark earlier in the document.
# ------------------------------------------------------------------
$goBackOld = $d.Bookmarks("_GoBack")
$goBackOld.Delete()

$itsRange = $d.Content
$itsRange.Find.Execute("it’s quite amazing")
$itsEnd = $itsRange.Start + 4   # "it’s" is 4 characters long
$goBackSpot = $d.Range($itsEnd, $itsEnd)
$d.Bookmarks.Add("_GoBack", $goBackSpot)
